# Saldo.xlsx update — refresh account balances in the "Export" sheet.
#
# The sheet is a flat list of (Conta, Nome, Saldo) rows sorted by
# descending Saldo. This script:
#   1. Refreshes one row in place (account 005009947 -> 008004995, new
#      name/balance).
#   2. Inserts two new rows for accounts that now have a balance putting
#      them higher in the sorted list (005341184, 005324981).
#   3. Removes the now-stale rows for accounts whose old balances no
#      longer apply (the old small-value 008004995/005324981 rows plus
#      several accounts that dropped out of the export entirely).
#
# Row numbers below are the original (pre-edit) 1-based sheet rows,
# applied bottom-to-top for the structural (insert/delete) operations so
# earlier row numbers stay valid while the row count shifts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-AccountRow($row, $conta, $nome, $saldo) {
    # Force column A to Text first so the zero-padded account number
    # ("008004995", ...) keeps its leading zeros instead of being
    # reinterpreted as a number.
    $ws.Range("A$row").NumberFormat = "@"
    $ws.Range("A$row").Value = $conta
    $ws.Range("B$row").Value = $nome
    $ws.Range("C$row").Value = $saldo
}

# --- Deletions (bottom-to-top so higher row numbers don't move first) ---
$ws.Rows(377).Delete()   # 004361159 / HFR       / -38178.76
$ws.Rows(376).Delete()   # 004352384 / BRASFORT  / -16646.54
$ws.Rows(374).Delete()   # 004415557 / FILIPE    / -5128.08
$ws.Rows(373).Delete()   # 004935287 / ODILON    / -5054.22
$ws.Rows(371).Delete()   # 005324981 / JO        / -1215.07  (stale; replaced below)
$ws.Rows(334).Delete()   # 008004995 / JOSE      / 3.27      (stale; replaced below)
$ws.Rows(328).Delete()   # 004346716 / TIAGO     / 4.8

# --- Insertions ---
$ws.Rows(68).Insert()    # new row for 005324981, pushes 004212438 down
Set-AccountRow 68 "005324981" "JO" 3758.51

$ws.Rows(51).Insert()    # new row for 005341184, pushes 004565108 down
Set-AccountRow 51 "005341184" "BRENO" 5312.64

# --- In-place update ---
Set-AccountRow 5 "008004995" "JOSE" 63753.27
